# Apply the edit described by the diff:
#  - Add a new shared string "pointOfSale"
#  - On the "InvalidLogin" sheet, append a new row 3: A3="admin", B3="pointOfSale"
#  - Autofit column B width on "InvalidLogin"
#  - Update the active-cell selection on "ValidLogin" (B3 -> A3) and "InvalidLogin" (B2 -> A3)

$wb = $excel.ActiveWorkbook

$wsValid = $wb.Worksheets.Item("ValidLogin")
$wsInvalid = $wb.Worksheets.Item("InvalidLogin")

# Add the new data row on the InvalidLogin sheet
$wsInvalid.Range("A3").Value = "admin"
$wsInvalid.Range("B3").Value = "pointOfSale"

# Autofit column B so its width reflects the new, wider content
$wsInvalid.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Update the recorded selection on ValidLogin (B3 -> A3)
$wsValid.Activate() | Out-Null
$wsValid.Range("A3").Select() | Out-Null

# Update the recorded selection on InvalidLogin (B2 -> A3) and leave it
# as the active sheet/tab, matching the workbook's original active tab
$wsInvalid.Activate() | Out-Null
$wsInvalid.Range("A3").Select() | Out-Null
